$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal number-looking string;
# force them to remain text (matching the source datas text-typed Price column)
# by pre-setting the cell number format to Text before assigning the value.
$textPriceRows = 5,6,8,11,14,15,16,19,20,21,22,23,24,25,26,27,31,34,37,38,39,40,41,45,46,47,48,49,51
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "97.973.67"
$ws.Range("E2").Value = "  +4.86%  "
$ws.Range("D3").Value = "3.142.28"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "241.39"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "609.71"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "0.382"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "3.139.21"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").Value = "0.788"
$ws.Range("E11").Value = "  -4.46%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "97.549.11"
$ws.Range("E13").Value = "  +4.76%  "
$ws.Range("D14").Value = "0.0000240"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "33.96"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "5.41"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "3.726.42"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "3.140.63"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "521.75"
$ws.Range("E19").Value = "  +18.36%  "
$ws.Range("D20").Value = "3.43"
$ws.Range("E20").Value = "  -6.75%  "
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "5.71"
$ws.Range("E22").Value = "  -4.78%  "
$ws.Range("D23").Value = "0.0000192"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("D24").Value = "8.79"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "88.51"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").Value = "5.46"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").Value = "11.62"
$ws.Range("E27").Value = "  -9.58%  "
$ws.Range("D28").Value = "3.313.83"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("E30").Value = "  -6.56%  "
$ws.Range("D31").Value = "0.176"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "8.95"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("D37").Value = "7.22"
$ws.Range("E37").Value = "  -9.00%  "
$ws.Range("D38").Value = "24.35"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "0.435"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("D41").Value = "466.09"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("E42").Value = "  -5.40%  "
$ws.Range("E43").Value = "  -11.38%  "
$ws.Range("D45").Value = "3.12"
$ws.Range("E45").Value = "  -4.67%  "
$ws.Range("D46").Value = "162.13"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").Value = "  +3.80%  "
$ws.Range("D48").Value = "0.695"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").Value = "4.51"
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "0.780"
$ws.Range("E51").Value = "  +7.13%  "
